$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "Euclides-Mecanica material"
$ws.Range("E4").Value = "João Rodrigues-CAD"
$ws.Range("B6").Value = "-"
$ws.Range("F6").Value = "Pedro Francisco-MTRM"
$ws.Range("F7").Value = "Pedro Francisco-MTRM"
